# Generate Report for Handback
#
# The handback for the "3e39e2b5-076c-4dd0-b50d-dc369fb46b7a" source file
# failed transform (file name mismatch), so the report is updated to show
# the failure status and the error detail for both the zh-cn and de-de
# target languages (and the rolled-up Overview sheet).

$wb = $excel.ActiveWorkbook

$ws_overview = $wb.Worksheets.Item("Overview")
$ws_zhcn     = $wb.Worksheets.Item("zh-cn")
$ws_dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"

$zhErrorDetail = "Handback file name: znd0xwua.nay is different with handoff file name: 3e39e2b5-076c-4dd0-b50d-dc369fb46b7a.31e0f3b0b1329d0b473ca3f8ea5bc23fca22cbcf.zh-cn."
$deErrorDetail = "Handback file name: znd0xwua.nay is different with handoff file name: 3e39e2b5-076c-4dd0-b50d-dc369fb46b7a.31e0f3b0b1329d0b473ca3f8ea5bc23fca22cbcf.de-de."

# Overview sheet: row 3 is the 3e39e2b5-... file, columns E (zh-cn) and F (de-de)
# hold that file's per-language status.
$ws_overview.Range("E3").Value = $newStatus
$ws_overview.Range("F3").Value = $newStatus

# zh-cn detail sheet: row 3 Status column (C) + new Error Detail (P).
$ws_zhcn.Range("C3").Value = $newStatus
$ws_zhcn.Range("P3").Value = $zhErrorDetail
$ws_zhcn.Columns.Item(16).ColumnWidth = 39.17

# de-de detail sheet: row 3 Status column (C) + new Error Detail (P).
$ws_dede.Range("C3").Value = $newStatus
$ws_dede.Range("P3").Value = $deErrorDetail
$ws_dede.Columns.Item(16).ColumnWidth = 39.17
